$wb = $excel.ActiveWorkbook

# Row 10 (ALC): A Jawbreaking Weapon of Staggering Weight
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 9999
$ws.Range("J10").Value = 9999
$ws.Range("L10").Value = 9999
$ws.Range("N10").Value = -10585

# Row 19 (ALC): Unbreak My Heart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 985.25
$ws.Range("I19").Value = 980.6667
$ws.Range("J19").Value = 999
$ws.Range("K19").Value = 980.6667
$ws.Range("L19").Value = 999
$ws.Range("M19").Value = -805.6667
$ws.Range("N19").Value = -1349

# Row 64 (ALC): Forged from the Void
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null

# Row 67 (ALC): Dodging the Draft (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null

# Row 70 (ALC): Consecrating Congregation
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11426.714
$ws.Range("I70").Value = 6697.3335
$ws.Range("K70").Value = 20092.0005
$ws.Range("M70").Value = -19822.0005

# Row 73 (ALC): Curbing the Contagion (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 11426.714
$ws.Range("I73").Value = 6697.3335
$ws.Range("K73").Value = 20092.0005
$ws.Range("M73").Value = -19156.0005

# Row 98 (ALC): The Dotted Line
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2882.0625
$ws.Range("I98").Value = 901.1111
$ws.Range("J98").Value = 5429
$ws.Range("K98").Value = 901.1111
$ws.Range("L98").Value = 5429
$ws.Range("M98").Value = 596.8889
$ws.Range("N98").Value = -8425

# Row 112 (ALC): Making Ends Meet
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 699
$ws.Range("I112").Value = 699
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 2097
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -989
$ws.Range("N112").Value = $null

# Row 113 (ALC): Amaro Kart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2933.4443
$ws.Range("I113").Value = 2561.2856
$ws.Range("J113").Value = 4236
$ws.Range("K113").Value = 2561.2856
$ws.Range("L113").Value = 4236
$ws.Range("M113").Value = 692.7143999999998
$ws.Range("N113").Value = -10744

# Row 122 (ALC): Wishful Inking
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2882.0625
$ws.Range("I122").Value = 901.1111
$ws.Range("J122").Value = 5429
$ws.Range("K122").Value = 2703.3333
$ws.Range("L122").Value = 16287
$ws.Range("M122").Value = -253.3332999999998
$ws.Range("N122").Value = -21187

# Row 2 (ARM): Ain't Got No Ingots
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 546.6667
$ws.Range("I2").Value = 369.0909
$ws.Range("K2").Value = 369.0909
$ws.Range("M2").Value = -256.0909

# Row 5 (ARM): The Alloyed Truth
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 447.5
$ws.Range("I5").Value = 516.8
$ws.Range("J5").Value = 101
$ws.Range("K5").Value = 516.8
$ws.Range("L5").Value = 101
$ws.Range("M5").Value = -404.8
$ws.Range("N5").Value = -325

# Row 32 (ARM): Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4855.8887
$ws.Range("I32").Value = 4855.8887
$ws.Range("K32").Value = 4855.8887
$ws.Range("M32").Value = -4568.8887

# Row 46 (ARM): Get Me the Usual
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 10049
$ws.Range("J46").Value = 10049
$ws.Range("L46").Value = 10049
$ws.Range("N46").Value = -10687

# Row 61 (ARM): Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4994.1113
$ws.Range("I61").Value = 1550.25
$ws.Range("K61").Value = 1550.25
$ws.Range("M61").Value = -1338.25

# Row 97 (ARM): Ore for Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 843.8
$ws.Range("I97").Value = 985
$ws.Range("J97").Value = 749.6667
$ws.Range("K97").Value = 985
$ws.Range("L97").Value = 749.6667
$ws.Range("M97").Value = -489
$ws.Range("N97").Value = -1741.6667

# Row 110 (ARM): Scheduled Maintenance
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1704.75
$ws.Range("I110").Value = 1704.75
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1704.75
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 340.25
$ws.Range("N110").Value = $null

# Row 116 (ARM): No Scope
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 546.6667
$ws.Range("I116").Value = 369.0909
$ws.Range("K116").Value = 369.0909
$ws.Range("M116").Value = 1924.9091

# Row 124 (ARM): Ace of Gloves
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 140000
$ws.Range("J124").Value = 140000
$ws.Range("L124").Value = 140000
$ws.Range("N124").Value = -149820

# Row 136 (ARM): Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4994.1113
$ws.Range("I136").Value = 1550.25
$ws.Range("K136").Value = 4650.75
$ws.Range("M136").Value = -2100.75

# Row 3 (BSM): Hells Bells
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 546.6667
$ws.Range("I3").Value = 369.0909
$ws.Range("K3").Value = 369.0909
$ws.Range("M3").Value = -255.0909

# Row 4 (BSM): Mending Fences
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 447.5
$ws.Range("I4").Value = 516.8
$ws.Range("J4").Value = 101
$ws.Range("K4").Value = 516.8
$ws.Range("L4").Value = 101
$ws.Range("M4").Value = -401.8
$ws.Range("N4").Value = -331

# Row 22 (BSM): Riveting Run
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 892
$ws.Range("I22").Value = 784.5
$ws.Range("K22").Value = 784.5
$ws.Range("M22").Value = -611.5

# Row 105 (BSM): Ingot to Wing It
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1729.2
$ws.Range("I105").Value = 1411.75
$ws.Range("K105").Value = 1411.75
$ws.Range("M105").Value = 335.25

# Row 139 (BSM): Maul Me
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

# Row 132 (CRP): Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1399.7368
$ws.Range("I132").Value = 1258.697
$ws.Range("K132").Value = 3776.090999999999
$ws.Range("M132").Value = -1246.090999999999

# Row 2 (CUL): Pork Is a Salty Food
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 201118.1
$ws.Range("I2").Value = 183433.33
$ws.Range("J2").Value = 222339.8
$ws.Range("K2").Value = 1100599.98
$ws.Range("L2").Value = 1334038.8
$ws.Range("M2").Value = -1100486.98
$ws.Range("N2").Value = -1334264.8

# Row 55 (CUL): Pagan Pastries
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = $null

# Row 130 (CUL): Blast from the Pasta
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2588.25
$ws.Range("I130").Value = 1937.5
$ws.Range("J130").Value = 3239
$ws.Range("K130").Value = 5812.5
$ws.Range("L130").Value = 9717
$ws.Range("M130").Value = -792.5
$ws.Range("N130").Value = -19757

# Row 134 (CUL): Don't Knock It Till You've Tried It
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 12841.889
$ws.Range("I134").Value = 1931.6666
$ws.Range("K134").Value = 5794.9998
$ws.Range("M134").Value = -724.9997999999996

# Row 137 (CUL): Creative Chocolate
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2204.389
$ws.Range("J137").Value = 2648.111
$ws.Range("L137").Value = 7944.333
$ws.Range("N137").Value = -18144.333

# Row 2 (GSM): Copper and Robbers
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 290.6
$ws.Range("I2").Value = 357.125
$ws.Range("J2").Value = 24.5
$ws.Range("K2").Value = 357.125
$ws.Range("L2").Value = 24.5
$ws.Range("M2").Value = -244.125
$ws.Range("N2").Value = -250.5

# Row 102 (GSM): Put the Metal to the Peddle
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1105.1538
$ws.Range("I102").Value = 991.5454999999999
$ws.Range("J102").Value = 1730
$ws.Range("K102").Value = 991.5454999999999
$ws.Range("L102").Value = 1730
$ws.Range("M102").Value = 630.4545000000001
$ws.Range("N102").Value = -4974

# Row 122 (GSM): Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3720.5715
$ws.Range("I122").Value = 2226.4614
$ws.Range("J122").Value = 6148.5
$ws.Range("K122").Value = 6679.3842
$ws.Range("L122").Value = 18445.5
$ws.Range("M122").Value = -4229.3842
$ws.Range("N122").Value = -23345.5

# Row 7 (LTW): Tan Before the Ban
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7190.76
$ws.Range("I7").Value = 2396.8
$ws.Range("K7").Value = 2396.8
$ws.Range("M7").Value = -2284.8

# Row 22 (LTW): Skin off Their Backs
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1412
$ws.Range("J22").Value = 1499.25
$ws.Range("L22").Value = 1499.25
$ws.Range("N22").Value = -2089.25

# Row 27 (LTW): Fire and Hide
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1412
$ws.Range("J27").Value = 1499.25
$ws.Range("L27").Value = 1499.25
$ws.Range("N27").Value = -1713.25

# Row 40 (LTW): Best Served Toad
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4148.3335
$ws.Range("J40").Value = 4322.5
$ws.Range("L40").Value = 4322.5
$ws.Range("N40").Value = -4594.5

# Row 46 (LTW): Supply Side Logic
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 38071.145
$ws.Range("I46").Value = 167333
$ws.Range("K46").Value = 167333
$ws.Range("M46").Value = -167145

# Row 68 (LTW): You Could Say It's a Moving Target
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2752.25
$ws.Range("J68").Value = 3003
$ws.Range("L68").Value = 3003
$ws.Range("N68").Value = -4501

# Row 71 (LTW): They Call It Bloody Mary (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2752.25
$ws.Range("J71").Value = 3003
$ws.Range("L71").Value = 15015
$ws.Range("N71").Value = -22503

# Row 93 (LTW): Hide to Go Seek
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = $null

# Row 126 (LTW): Battered Books
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7190.76
$ws.Range("I126").Value = 2396.8
$ws.Range("K126").Value = 7190.400000000001
$ws.Range("M126").Value = -4720.400000000001

# Row 132 (LTW): Tenets of Tanning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3206.35
$ws.Range("I132").Value = 2811.4666
$ws.Range("J132").Value = 4391
$ws.Range("K132").Value = 8434.399800000001
$ws.Range("L132").Value = 13173
$ws.Range("M132").Value = -5904.399800000001
$ws.Range("N132").Value = -18233

# Row 62 (WVR): Pride Up in Smoke
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14998.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null

# Row 65 (WVR): Desperate for Diversionaries (L)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 14998.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null

# Row 100 (WVR): Of Great Import
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 12503543
$ws.Range("I100").Value = 14288335
$ws.Range("K100").Value = 28576670
$ws.Range("M100").Value = -28576129

# Row 122 (WVR): Heavy Armoire
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 311.57144
$ws.Range("I122").Value = 311.57144
$ws.Range("K122").Value = 934.71432
$ws.Range("M122").Value = 1515.28568

# Row 132 (WVR): Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 991.8333
$ws.Range("I132").Value = 991.8333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2975.4999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -445.4998999999998
$ws.Range("N132").Value = $null
